# Applies the "Add final Code Changes" edit described by the OOXML diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Mark the five inline picture runs as NoProofing -> adds
#    <w:rPr><w:noProof/></w:rPr> to the run that hosts each <w:drawing>.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.NoProofing = 1
    }
}

# ---------------------------------------------------------------------
# 2) "Fibonacci – fib.j" -> "Fibonacci – " + proofErr(fib.j) split
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Fibonacci") | Out-Null
$p = $d.Paragraphs(18)
$xml = '<w:p w14:paraId="0C26D569" w14:textId="13E98443" w:rsidR="00894509" w:rsidRDefault="00894509" w:rsidP="00894509"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Fibonacci – </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>fib.j</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 3) "Factorial – fact.j" -> "Factorial – " + proofErr(fact.j) split
# ---------------------------------------------------------------------
$p = $d.Paragraphs(19)
$xml = '<w:p w14:paraId="0083A994" w14:textId="0113005D" w:rsidR="00894509" w:rsidRDefault="00894509" w:rsidP="00894509"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Factorial – </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>fact.j</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 4) "I added the for loop by extending the compile_stmt, and adding
#    support for For loops in Stmt parser." split around "stmt, and"
# ---------------------------------------------------------------------
$p = $d.Paragraphs(22)
$xml = '<w:p w14:paraId="51A0CDD8" w14:textId="066A69D0" w:rsidR="003D79DE" w:rsidRDefault="003D79DE" w:rsidP="003D79DE"><w:r><w:t>I added the for loop by extending the compile_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>stmt, and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> adding support for For loops in Stmt parser.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 5) Insert a new paragraph "The assembler instructions can be found in
#    nestedi.j" (bold+underline) right before "Currently since the
#    variables are not scoped..."
# ---------------------------------------------------------------------
$p = $d.Paragraphs(35)
$p.Range.InsertParagraphAfter()
$newp = $d.Paragraphs(36)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The assembler instructions can be found in </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>nestedi.j</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$newp.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 6) "We iterate through the inner loop until the condition is not met
#    which is when i=10" -> split with an inline OMath "i ≤ 10"
# ---------------------------------------------------------------------
$p = $d.Paragraphs(40)
$xml = '<w:p w14:paraId="0C3957B9" w14:textId="5EF2EAD2" w:rsidR="003D79DE" w:rsidRDefault="003D79DE" w:rsidP="003D79DE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">We iterate through the inner loop until the condition </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>i</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>&#8804;</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>10</m:t></m:r></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>is not met</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 7) "When we exit the inner loop, we compare this new I with the
#    condition which will also not be satisfied." -> isolate "i" as its
#    own run
# ---------------------------------------------------------------------
$p = $d.Paragraphs(41)
$xml = '<w:p w14:paraId="69459D08" w14:textId="77777777" w:rsidR="00B34C0A" w:rsidRDefault="00B34C0A" w:rsidP="003D79DE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">When we exit the inner loop, we compare this new </w:t></w:r><w:r><w:t>i</w:t></w:r><w:r><w:t xml:space="preserve"> with the condition which will also not be satisfied.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------
# 8) "Therefore the output will look like the following" -> split into
#    "Therefore" + "," + " the output will look like the following"
# ---------------------------------------------------------------------
$p = $d.Paragraphs(43)
$xml = '<w:p w14:paraId="19FF41F3" w14:textId="476B11DB" w:rsidR="00B34C0A" w:rsidRDefault="00B34C0A" w:rsidP="003D79DE"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Therefore</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> the output will look like the following</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

Write-Output "done"
